# "Generate Report for Handback" - refresh the timestamps/status recorded by
# the handback status report, mirroring what the generator would produce on
# a fresh run (a few minutes later than the previous run, and the zh-cn
# priority flipped from "ht" to "mt").

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
# "Latest HO Xliff Generate Date" for 029a50fe-...md
$wsOverview.Range("G2").Value = "2016-08-18 00:15:21"
# "Latest HO Xliff Generate Date" for af3e64a7-...md
$wsOverview.Range("G4").Value = "2016-08-18 00:15:21"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority: ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-18 00:15:16"
$wsZhCn.Range("H4").Value = "2016-08-18 00:15:16"
# Correspond Handback DateTime
$wsZhCn.Range("K2").Value = "2016-08-18 00:15:34"
$wsZhCn.Range("K4").Value = "2016-08-18 00:15:34"

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority: ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
# Correspond Handoff Datetime
$wsDeDe.Range("H2").Value = "2016-08-18 00:15:21"
$wsDeDe.Range("H4").Value = "2016-08-18 00:15:21"
# Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-18 00:15:41"
$wsDeDe.Range("K4").Value = "2016-08-18 00:15:41"
